$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 123: Nearly Bare / Gaja Grimoire
$ws.Range("H123").Value = 41986.668
$ws.Range("J123").Value = 41986.668
$ws.Range("L123").Value = 41986.668
$ws.Range("N123").Value = -51786.668

# Row 127: Liquid Competence / Competent Craftsman's Draught
$ws.Range("H127").Value = 2222.5
$ws.Range("I127").Value = 889.1
$ws.Range("J127").Value = 2639.1875
$ws.Range("K127").Value = 2667.3
$ws.Range("L127").Value = 7917.5625
$ws.Range("M127").Value = 2292.7
$ws.Range("N127").Value = -17837.5625

# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 1893.0615
$ws.Range("I132").Value = 1765.1864
$ws.Range("J132").Value = 3150.5
$ws.Range("K132").Value = 5295.5592
$ws.Range("L132").Value = 9451.5
$ws.Range("M132").Value = -2765.5592
$ws.Range("N132").Value = -14511.5

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 4001.2122
$ws.Range("I138").Value = 875
$ws.Range("J138").Value = 4743.6875
$ws.Range("K138").Value = 2625
$ws.Range("L138").Value = 14231.0625
$ws.Range("M138").Value = 2515
$ws.Range("N138").Value = -24511.0625

# Row 141: Remedy for Reason / Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 3805
$ws.Range("I141").Value = 707.5
$ws.Range("J141").Value = 10000
$ws.Range("K141").Value = 2122.5
$ws.Range("L141").Value = 30000
$ws.Range("M141").Value = 3057.5
$ws.Range("N141").Value = -40360

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 5957778
$ws.Range("I32").Value = 6807706
$ws.Range("K32").Value = 6807706
$ws.Range("M32").Value = -6807419

# Row 138: Don't Ask about the Rivets / Titanium Gold Helm of Casting
$ws.Range("H138").Value = 41429
$ws.Range("J138").Value = 41429
$ws.Range("L138").Value = 41429
$ws.Range("N138").Value = -51709

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 2044.7142
$ws.Range("I31").Value = 2207.2432
$ws.Range("J31").Value = 1543.5834
$ws.Range("K31").Value = 2207.2432
$ws.Range("L31").Value = 1543.5834
$ws.Range("M31").Value = -1912.2432
$ws.Range("N31").Value = -2133.5834

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 2044.7142
$ws.Range("I34").Value = 2207.2432
$ws.Range("J34").Value = 1543.5834
$ws.Range("K34").Value = 2207.2432
$ws.Range("L34").Value = 1543.5834
$ws.Range("M34").Value = -2005.2432
$ws.Range("N34").Value = -1947.5834

# Row 55: Ready for a Rematch / Mythril Lance
$ws.Range("H55").Value = 8000
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 8000
$ws.Range("K55").Value = 0
$ws.Range("L55").ClearContents()
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -8630

# Row 94: Beech, Please / Beech Lumber
$ws.Range("H94").Value = 4026.6667
$ws.Range("I94").Value = 710
$ws.Range("K94").Value = 710
$ws.Range("M94").Value = -259

$ws = $wb.Worksheets.Item("CUL")
# Row 11: Putting the Squeeze On / Orange Juice
$ws.Range("H11").Value = 50070
$ws.Range("I11").Value = 66726.664
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 200179.992
$ws.Range("L11").Value = 300
$ws.Range("M11").Value = -200039.992
$ws.Range("N11").Value = -580

# Row 51: The Perks of Life at Sea / Jerked Beef
$ws.Range("H51").Value = 1535.3334
$ws.Range("I51").Value = 242.4
$ws.Range("J51").Value = 8000
$ws.Range("K51").Value = 727.2
$ws.Range("L51").Value = 24000
$ws.Range("M51").Value = -267.2
$ws.Range("N51").Value = -24920

# Row 68: Such a Butter Face / Fermented Butter
$ws.Range("H68").Value = 880
$ws.Range("J68").Value = 1000
$ws.Range("L68").Value = 3000
$ws.Range("N68").Value = -4622

# Row 71: No Margarine of Error (L) / Fermented Butter
$ws.Range("H71").Value = 880
$ws.Range("J71").Value = 1000
$ws.Range("L71").Value = 9000
$ws.Range("N71").Value = -17112

# Row 80: Saucy for a Suitor / Hollandaise Sauce
$ws.Range("H80").Value = 8000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 8000
$ws.Range("K80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -25872

# Row 83: Saved by the Sauce (L) / Hollandaise Sauce
$ws.Range("H83").Value = 8000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 8000
$ws.Range("K83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -81360

# Row 107: Slippery Service / Frantoio Oil
$ws.Range("H107").Value = 486913.7
$ws.Range("I107").Value = 1010.8889
$ws.Range("J107").Value = 1111645.9
$ws.Range("K107").Value = 3032.6667
$ws.Range("L107").Value = 3334937.7
$ws.Range("M107").Value = -1112.6667
$ws.Range("N107").Value = -3338777.7

# Row 132: More Mezcal / Cooking Mezcal
$ws.Range("H132").Value = 926.8570999999999
$ws.Range("J132").Value = 1066.6666
$ws.Range("L132").Value = 9599.999400000001
$ws.Range("N132").Value = -14659.9994

# Row 136: Simple Is Hardest / Spaghetti al Olio e Peperoncino
$ws.Range("H136").Value = 4130
$ws.Range("I136").Value = 3092
$ws.Range("J136").Value = 4995
$ws.Range("K136").Value = 9276
$ws.Range("L136").Value = 14985
$ws.Range("M136").Value = -4176
$ws.Range("N136").Value = -25185

# Row 139: Najoothie / Wild Banana Blend
$ws.Range("H139").Value = 1791.4286
$ws.Range("I139").Value = 1791.4286
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 5374.2858
$ws.Range("L139").Value = 0
$ws.Range("M139").ClearContents()
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit / Mythrite Ingot
$ws.Range("H70").Value = 4844.4546
$ws.Range("I70").Value = 4048.3333
$ws.Range("J70").Value = 5799.8
$ws.Range("K70").Value = 4048.3333
$ws.Range("L70").Value = 5799.8
$ws.Range("M70").Value = -3778.3333
$ws.Range("N70").Value = -6339.8

# Row 73: Hulls of Broken Dreams (L) / Mythrite Ingot
$ws.Range("H73").Value = 4844.4546
$ws.Range("I73").Value = 4048.3333
$ws.Range("J73").Value = 5799.8
$ws.Range("K73").Value = 4048.3333
$ws.Range("L73").Value = 5799.8
$ws.Range("M73").Value = -3112.3333
$ws.Range("N73").Value = -7671.8

$ws = $wb.Worksheets.Item("LTW")
# Row 5: These Boots Are Made for Wailing / Leather Duckbills of Gathering
$ws.Range("H5").Value = 20000
$ws.Range("J5").Value = 20000
$ws.Range("L5").Value = 20000
$ws.Range("N5").Value = -20226

# Row 46: Supply Side Logic / Boar Leather
$ws.Range("H46").Value = 1187.2
$ws.Range("I46").Value = 1298
$ws.Range("J46").Value = 965.6
$ws.Range("K46").Value = 1298
$ws.Range("L46").Value = 965.6
$ws.Range("M46").Value = -1110
$ws.Range("N46").Value = -1341.6

# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 7022.3477
$ws.Range("I122").Value = 8088.7646
$ws.Range("J122").Value = 4000.8333
$ws.Range("K122").Value = 24266.2938
$ws.Range("L122").Value = 12002.4999
$ws.Range("M122").Value = -21816.2938
$ws.Range("N122").Value = -16902.4999

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 2198.842
$ws.Range("I132").Value = 1470.4546
$ws.Range("K132").Value = 4411.3638
$ws.Range("M132").Value = -1881.3638

# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 4414.3125
$ws.Range("I136").Value = 2012.4
$ws.Range("J136").Value = 8417.5
$ws.Range("K136").Value = 6037.200000000001
$ws.Range("L136").Value = 25252.5
$ws.Range("M136").Value = -3487.200000000001
$ws.Range("N136").Value = -30352.5

$ws = $wb.Worksheets.Item("WVR")
# Row 122: Heavy Armoire / Dark Hempen Cloth
$ws.Range("H122").Value = 2116.85
$ws.Range("I122").Value = 1938
$ws.Range("K122").Value = 5814
$ws.Range("M122").Value = -3364

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 1566.7963
$ws.Range("I136").Value = 1681.8158
$ws.Range("J136").Value = 1293.625
$ws.Range("K136").Value = 5045.4474
$ws.Range("L136").Value = 3880.875
$ws.Range("M136").Value = -2495.4474
$ws.Range("N136").Value = -8980.875

# Row 141: Silk for Sunperch / Thunderyards Silk Coat of Casting
$ws.Range("H141").Value = 73282.14
$ws.Range("J141").Value = 73282.14
$ws.Range("L141").Value = 73282.14
$ws.Range("N141").Value = -83642.14
